# Work_Progress.xlsx update
#  - Add new rows of activity log to Sheet1 (rows 20-29)
#  - Add a new "Sheet2" worksheet capturing the "Bugs" note
#  - Update selections to match the final saved state

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: fill in the missing dates for the existing "testing the screen" row ---
$ws1.Range("B20").Value = 45248
$ws1.Range("B20").NumberFormat = $ws1.Range("B19").NumberFormat
$ws1.Range("C20").Value = 45248
$ws1.Range("C20").NumberFormat = $ws1.Range("C19").NumberFormat

# --- Sheet1: append the new activity rows (order chosen to match the
#     original shared-string insertion order of the authored workbook) ---
$ws1.Range("A22").Value = "Create Employee Screen Client side "
$ws1.Range("A24").Value = "Create Employee Screen Validations"
$ws1.Range("A23").Value = "Create Employee Screen Servier side "
$ws1.Range("A25").Value = "Modify Employee "
$ws1.Range("A26").Value = "Navigation from List to Create screen"
$ws1.Range("A21").Value = "demo the screen "
$ws1.Range("A29").Value = "split the utility into a separate application "

# --- Add Sheet2 right after Sheet1 ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "Bugs"
$ws2.Range("A2").Value = "no error thrown but said import success"
$ws2.Range("B2").Value = "fixed"

# These two reuse strings introduced above / already present in the sheet
$ws1.Range("A27").Value = "testing the screen "
$ws1.Range("A28").Value = "demo the screen "

# Column A on Sheet2 needs to be wide enough to show the bug description
$ws2.Columns.Item(1).ColumnWidth = 35.9

# --- Restore selections to match the final workbook state ---
$ws2.Range("B3").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("D20").Select() | Out-Null
